$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue 'D2' '69.905.65'
$ws.Range('E2').Value = '  -1.89%  '
Set-TextValue 'D3' '3.709.97'
$ws.Range('E3').Value = '  -2.17%  '
$ws.Range('E4').Value = '  +0.06%  '
Set-TextValue 'D5' '618.39'
$ws.Range('E5').Value = '  -0.07%  '
Set-TextValue 'D6' '179.06'
$ws.Range('E6').Value = '  -0.97%  '
Set-TextValue 'D7' '3.711.68'
$ws.Range('E7').Value = '  -1.69%  '
$ws.Range('E8').Value = '  -0.14%  '
Set-TextValue 'D9' '0.532'
$ws.Range('E9').Value = '  -2.05%  '
$ws.Range('E10').Value = '  +0.67%  '
Set-TextValue 'D11' '6.29'
$ws.Range('E11').Value = '  -3.72%  '
Set-TextValue 'D12' '0.483'
$ws.Range('E12').Value = '  -3.73%  '
Set-TextValue 'D13' '40.41'
$ws.Range('E13').Value = '  -0.01%  '
Set-TextValue 'D14' '0.0000256'
$ws.Range('E14').Value = '  -0.18%  '
Set-TextValue 'D15' '4.345.30'
$ws.Range('E15').Value = '  -1.92%  '
Set-TextValue 'D16' '3.719.93'
$ws.Range('E16').Value = '  -2.33%  '
Set-TextValue 'D17' '69.939.40'
$ws.Range('E17').Value = '  -2.04%  '
$ws.Range('E18').Value = '  -1.88%  '
Set-TextValue 'D19' '7.58'
$ws.Range('E19').Value = '  +0.18%  '
Set-TextValue 'D20' '16.60'
$ws.Range('E20').Value = '  -1.89%  '
Set-TextValue 'D21' '503.30'
$ws.Range('E21').Value = '  -3.35%  '
Set-TextValue 'D22' '9.25'
$ws.Range('E22').Value = '  -1.62%  '
Set-TextValue 'D23' '0.718'
$ws.Range('E23').Value = '  -4.36%  '
$ws.Range('E24').Value = '  -0.14%  '
Set-TextValue 'D25' '86.31'
$ws.Range('E25').Value = '  -2.91%  '
Set-TextValue 'D26' '11.34'
$ws.Range('E26').Value = '  +1.95%  '
Set-TextValue 'D27' '13.03'
$ws.Range('E27').Value = '  -3.92%  '
Set-TextValue 'D28' '0.0000135'
$ws.Range('E28').Value = '  +19.41%  '
Set-TextValue 'D29' '0.997'
$ws.Range('E29').Value = '  -0.35%  '
Set-TextValue 'D30' '2.47'
$ws.Range('E30').Value = '  -2.95%  '
Set-TextValue 'D31' '2.92'
$ws.Range('E31').Value = '  +0.56%  '
Set-TextValue 'D32' '7.89'
$ws.Range('E32').Value = '  -2.64%  '
Set-TextValue 'D33' '30.96'
$ws.Range('E33').Value = '  -4.90%  '
Set-TextValue 'D34' '0.114'
$ws.Range('E34').Value = '  -1.72%  '
Set-TextValue 'D35' '1.00'
$ws.Range('E35').Value = '  +0.07%  '
Set-TextValue 'D36' '1.05'
$ws.Range('E36').Value = '  -0.39%  '
Set-TextValue 'D37' '6.12'
$ws.Range('E37').Value = '  -0.69%  '
Set-TextValue 'D38' '0.137'
$ws.Range('E38').Value = '  +1.79%  '
Set-TextValue 'D39' '0.337'
$ws.Range('E39').Value = '  -1.79%  '
Set-TextValue 'D40' '2.07'
$ws.Range('E40').Value = '  -7.58%  '
Set-TextValue 'D41' '49.97'
$ws.Range('E41').Value = '  -3.30%  '
Set-TextValue 'D42' '45.30'
$ws.Range('E42').Value = '  +2.03%  '
Set-TextValue 'D43' '429.36'
$ws.Range('E43').Value = '  -3.02%  '
Set-TextValue 'D44' '2.89'
$ws.Range('E44').Value = '  +2.56%  '
Set-TextValue 'D45' '8.65'
$ws.Range('E45').Value = '  -2.48%  '
Set-TextValue 'D46' '2.979.82'
$ws.Range('E46').Value = '  -6.32%  '
Set-TextValue 'D47' '0.0361'
$ws.Range('E47').Value = '  -1.49%  '
Set-TextValue 'D48' '27.38'
$ws.Range('E48').Value = '  -2.24%  '
$ws.Range('E49').Value = '  -0.04%  '
Set-TextValue 'D50' '136.83'
$ws.Range('E50').Value = '  -2.36%  '
Set-TextValue 'D51' '2.48'
$ws.Range('E51').Value = '  +0.63%  '
